$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Rename the "Hangman" minigame idea to "mastermind".
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Hangman", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "mastermind", 2)

# Locate that paragraph again (it is the "List Paragraph" that now reads
# "mastermind").
$pMastermind = $d.Paragraphs(22)
$pNext = $pMastermind.Next()

# ---------------------------------------------------------------------
# 2. Move the (hidden) _GoBack bookmark from the following empty
#    paragraph to the end of the "mastermind" run, collapsed, with
#    nothing selected - exactly like Word leaves it after the last edit.
#    A bookmark collapsed exactly on a paragraph mark gets mis-placed by
#    this host, so we briefly insert a throw-away character to give the
#    engine a real text position to collapse onto, then remove it again.
# ---------------------------------------------------------------------
$r = $pMastermind.Range
$insertPos = $r.Start + $r.Text.TrimEnd("`r").Length

$placeholder = $d.Range($insertPos, $insertPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($insertPos, $insertPos)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$d.Range($insertPos, $insertPos + 1).Delete() | Out-Null

# ---------------------------------------------------------------------
# 3. The paragraph right after "mastermind" (previously empty, holder of
#    the bookmark) now becomes the "Simon Say" idea. Seed it with a run
#    that carries the correct run formatting (lang=en-US) by cloning the
#    formatted text from the "mastermind" run, then swap in the new
#    wording.
# ---------------------------------------------------------------------
$pNext.Range.FormattedText = $pMastermind.Range.FormattedText
$pNext.Range.Text = "Simon Say"
$pSimon = $pNext

# ---------------------------------------------------------------------
# 4. Add a brand-new list item after it for "Breakout", inheriting the
#    same list formatting (pStyle / numPr / rPr) via InsertParagraphAfter.
# ---------------------------------------------------------------------
$pSimon.Range.InsertParagraphAfter()
$pBreakout = $pSimon.Next()
$pBreakout.Range.Text = "Breakout"

Write-Host "done"
